# daily auto push: 2026-01-07 13:49 UTC
#
# This sheet is a rolling daily log: every day a fresh row of readings is
# inserted right after the "today" block (row 587) and everything that was
# there before (2026/12/29 .. 2027/01/05) shifts down by one row, growing
# the used range from A1:D628 to A1:D629.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 587:628 down to 588:629 and open up a blank row 587.
$ws.Rows.Item(587).Insert()

# A587 holds a date formatted as plain text ("2026/01/07"), exactly like the
# surrounding rows (e.g. A586). Flip the cell to Text first so Excel doesn't
# silently turn the string into a date serial number, then drop the
# formatting override again so the cell keeps the sheet's plain default
# style (no explicit number format), matching the rest of the column.
$ws.Range("A587").NumberFormat = "@"
$ws.Range("A587").Value = "2026/01/07"
$ws.Range("A587").ClearFormats()

$ws.Range("B587").Value = "水"
$ws.Range("C587").Value = 20
$ws.Range("D587").Value = 201
